$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: 自动登录显示 cardnum ------------------------------------
$ws.Range("A1").Copy($ws.Range("A13"))
$ws.Range("A13").Value = "自动登录显示 cardnum"
$ws.Range("B8").Copy($ws.Range("B13"))
$ws.Range("B13").Value = 42792

# --- Row 14: 局部刷新按钮禁止 ------------------------------------------
$ws.Range("A1").Copy($ws.Range("A14"))
$ws.Range("A14").Value = "局部刷新按钮禁止"
$ws.Range("B8").Copy($ws.Range("B14"))
$ws.Range("B14").Value = 42792

# --- Row 15: 校园网已用流量 --------------------------------------------
$ws.Range("A1").Copy($ws.Range("A15"))
$ws.Range("A15").Value = "校园网已用流量"
$ws.Range("B8").Copy($ws.Range("B15"))
$ws.Range("B15").Value = 42792

# --- New "备注" (remarks) header column --------------------------------
$ws.Range("D1").Copy($ws.Range("E1"))
$ws.Range("E1").Value = "备注"

# --- Row 16: 人文讲座预告 (resolved) -----------------------------------
$ws.Range("A1").Copy($ws.Range("A16"))
$ws.Range("A16").Value = "人文讲座预告"
$ws.Range("B8").Copy($ws.Range("B16"))
$ws.Range("B16").Value = 42792
$ws.Range("C2").Copy($ws.Range("C16"))
$ws.Range("C16").Value = "√"
$ws.Range("B8").Copy($ws.Range("D16"))
$ws.Range("D16").Value = 42792

# --- Row 17: 修改密码情况 ----------------------------------------------
$ws.Range("A1").Copy($ws.Range("A17"))
$ws.Range("A17").Value = "修改密码情况"

# --- Row 18: 图书查询 (resolved) ---------------------------------------
$ws.Range("A1").Copy($ws.Range("A18"))
$ws.Range("A18").Value = "图书查询"
$ws.Range("B8").Copy($ws.Range("B18"))
$ws.Range("B18").Value = 42792
$ws.Range("C2").Copy($ws.Range("C18"))
$ws.Range("C18").Value = "√"
$ws.Range("B8").Copy($ws.Range("D18"))
$ws.Range("D18").Value = 42792

# --- Row heights to match the target layout -----------------------------
$ws.Rows("13").RowHeight = 28.5
$ws.Rows("14").RowHeight = 34.5
$ws.Rows("15").RowHeight = 42
$ws.Rows("16").RowHeight = 22.5
$ws.Rows("17").RowHeight = 23.25

# --- Update view state: scroll down and move the active selection -------
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("G20").Select()
